$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "67.235.72"
Set-TextValue "E2" "  +0.62%  "

Set-TextValue "D3" "3.483.59"
Set-TextValue "E3" "  -0.15%  "

Set-TextValue "D5" "593.74"
Set-TextValue "E5" "  +0.15%  "

Set-TextValue "D6" "178.69"
Set-TextValue "E6" "  +4.12%  "

Set-TextValue "E7" "  +0.01%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue "D8" "3.486.32"
Set-TextValue "E8" "  -0.05%  "

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D9" "0.595"
Set-TextValue "E9" "  +0.99%  "

Set-TextValue "D10" "0.138"
Set-TextValue "E10" "  +4.96%  "

Set-TextValue "E11" "  -2.33%  "

Set-TextValue "D12" "0.436"
Set-TextValue "E12" "  +1.07%  "

Set-TextValue "D13" "4.087.54"
Set-TextValue "E13" "  -0.14%  "

Set-TextValue "E14" "  +10.53%  "

Set-TextValue "E15" "  +1.35%  "

Set-TextValue "D16" "67.292.08"
Set-TextValue "E16" "  +0.68%  "

Set-TextValue "E17" "  -0.40%  "

Set-TextValue "D18" "3.485.42"
Set-TextValue "E18" "  -0.25%  "

Set-TextValue "D19" "6.26"
Set-TextValue "E19" "  -0.12%  "

Set-TextValue "D20" "14.30"

Set-TextValue "D21" "389.01"
Set-TextValue "E21" "  -0.95%  "

Set-TextValue "D22" "8.00"
Set-TextValue "E22" "  +0.51%  "

Set-TextValue "D23" "74.19"
Set-TextValue "E23" "  +1.69%  "

Set-TextValue "D24" "1.00"
Set-TextValue "E24" "  -0.02%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D25" "0.536"
Set-TextValue "E25" "  +0.34%  "

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D26" "5.72"
Set-TextValue "E26" "  +0.49%  "

Set-TextValue "D27" "0.0000121"
Set-TextValue "E27" "  +0.39%  "

Set-TextValue "D28" "10.37"
Set-TextValue "E28" "  +1.95%  "

Set-TextValue "E29" "  -3.47%  "

Set-TextValue "E30" "  +0.31%  "

Set-TextValue "D31" "6.15"
Set-TextValue "E31" "  -0.35%  "

Set-TextValue "D32" "1.42"
Set-TextValue "E32" "  -0.42%  "

Set-TextValue "E33" "  +0.44%  "

Set-TextValue "D34" "23.56"
Set-TextValue "E34" "  -0.44%  "

Set-TextValue "D35" "7.37"
Set-TextValue "E35" "  +0.65%  "

Set-TextValue "E36" "  +0.03%  "

Set-TextValue "E37" "  -1.60%  "

Set-TextValue "D38" "164.08"
Set-TextValue "E38" "  +0.89%  "

Set-TextValue "D39" "0.872"
Set-TextValue "E39" "  -0.55%  "

Set-TextValue "D40" "1.88"
Set-TextValue "E40" "  -0.72%  "

Set-TextValue "D41" "2.74"
Set-TextValue "E41" "  +7.41%  "

Set-TextValue "D42" "6.84"
Set-TextValue "E42" "  -1.82%  "

Set-TextValue "D43" "4.65"
Set-TextValue "E43" "  +0.28%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D44" "27.25"
Set-TextValue "E44" "  +0.57%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D45" "2.834.60"
Set-TextValue "E45" "  +0.68%  "

Set-TextValue "D46" "26.28"
Set-TextValue "E46" "  +0.54%  "

Set-TextValue "D47" "0.0722"
Set-TextValue "E47" "  -2.40%  "

Set-TextValue "D48" "41.64"
Set-TextValue "E48" "  -2.73%  "

Set-TextValue "D49" "0.0301"
Set-TextValue "E49" "  -0.38%  "

Set-TextValue "D50" "336.47"
Set-TextValue "E50" "  -0.01%  "

Set-TextValue "E51" "  -2.00%  "
